# Append: 2026-01-03 06:28 JST
# Re-scrape result: newest 5 rows kept (incl. header), older rows dropped,
# row contents shifted/refreshed, column widths tweaked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Remove obsolete rows 6-10 (sheet shrinks from H10 to H5) ---
$ws.Range("A6:H10").EntireRow.Delete()

# --- Refresh the 4 data rows (2-5) with the latest scrape values ---

# Row 2
$ws.Range("A2").Value = "2026-01-03 06:28:34"
$ws.Range("B2").Value = "【急募】ノーコードツール「bubble」でのメンテナンス依頼と次期開発依頼"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5464287"
$ws.Range("G2").Value = 135
$ws.Range("H2").Value = "◆ツール,開発"

# Row 3
$ws.Range("A3").Value = "2026-01-03 06:28:34"
$ws.Range("B3").Value = "Google Spread Sheetを用いた社内ツールの改修+α"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5464276"
$ws.Range("G3").Value = 98
$ws.Range("H3").Value = "◆ツール ◇社内ツール"

# Row 4
$ws.Range("A4").Value = "2026-01-03 06:28:34"
$ws.Range("B4").Value = "【外国人大歓迎】【急募】ECツールの保守・バグ修正・機能追加エンジニア募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5457026"
$ws.Range("G4").Value = 68
$ws.Range("H4").Value = "◆ツール"

# Row 5
$ws.Range("A5").Value = "2026-01-03 06:28:34"
$ws.Range("B5").Value = "【急募】Windows11対応ソフトウェア開発"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5464212"
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = "◆開発"

# --- Rebuild hyperlinks so only F2:F5 point at the (refreshed) URLs ---
# (Deleting hyperlinks through a Range clears the sheet's whole collection
# in this engine, so drop them all and re-add just the 4 we need.)
$ws.Range("F2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5464287")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5464276")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5457026")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5464212")

# Hyperlinks.Add() force-applies its own one-off "Hyperlink-ish" style (a
# duplicate cellXf) to the cell; reassign the plain named "Hyperlink" style
# on top so F2:F5 end up matching the original style index again.
$ws.Range("F2:F5").Style = "Hyperlink"

# --- Column width tweaks (values are in "characters"; COM's ColumnWidth
#     setter bakes in the usual +5/6 Maximum-Digit-Width padding, so back
#     it out to land on the exact stored widths from the target file) ---
$ws.Columns.Item(2).ColumnWidth = 229/6   # -> stored width 39
$ws.Columns.Item(4).ColumnWidth = 187/6   # -> stored width 32
$ws.Columns.Item(8).ColumnWidth = 73/6    # -> stored width 13
